# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header update ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 18:52"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 770981
$ws.Range("C4").Value = 6345
$ws.Range("E4").Value = 658143
$ws.Range("G4").Value = 774
$ws.Range("H4").Value = 41349

# --- Italia (row 6) ---
$ws.Range("B6").Value = 181228
$ws.Range("C6").Value = 2256
$ws.Range("D6").Value = 48877
$ws.Range("E6").Value = 108237
$ws.Range("F6").Value = 2573
$ws.Range("G6").Value = 454
$ws.Range("H6").Value = 24114

# --- Alemania (row 8) ---
$ws.Range("B8").Value = 146293
$ws.Range("C8").Value = 551
$ws.Range("E8").Value = 50110
$ws.Range("G8").Value = 41
$ws.Range("H8").Value = 4683

# --- Brasil (row 15) ---
$ws.Range("B15").Value = 39548
$ws.Range("C15").Value = 894
$ws.Range("E15").Value = 14911
$ws.Range("G15").Value = 45
$ws.Range("H15").Value = 2507

# --- Polonia / Ecuador swap (rows 30-31) ---
# Ecuador moves above Polonia with updated figures; Polonia keeps its
# previous figures but drops one rank.
$ws.Range("A30").Value = "Ecuador"
$ws.Range("B30").Value = 10128
$ws.Range("C30").Value = 660
$ws.Range("D30").Value = 1150
$ws.Range("E30").Value = 8471
$ws.Range("F30").Value = 124
$ws.Range("G30").Value = 33
$ws.Range("H30").Value = 507

$ws.Range("A31").Value = "Polonia"
$ws.Range("B31").Value = 9593
$ws.Range("C31").Value = 306
$ws.Range("D31").Value = 1133
$ws.Range("E31").Value = 8080
$ws.Range("F31").Value = 160
$ws.Range("G31").Value = 20
$ws.Range("H31").Value = 380

# --- Noruega (row 38) ---
$ws.Range("B38").Value = 7127
$ws.Range("C38").Value = 49
$ws.Range("E38").Value = 6914
$ws.Range("G38").Value = 16
$ws.Range("H38").Value = 181

# --- Kazajistan (row 66) ---
$ws.Range("B66").Value = 1852
$ws.Range("C66").Value = 176
$ws.Range("D66").Value = 447
$ws.Range("E66").Value = 1386

# --- Irak (row 69) ---
$ws.Range("B69").Value = 1574
$ws.Range("C69").Value = 35
$ws.Range("D69").Value = 1043
$ws.Range("E69").Value = 449

# --- Trinidad yTobago / Birmania swap (rows 137-138) ---
# Birmania moves above Trinidad yTobago with updated figures; Trinidad
# yTobago keeps its previous figures but drops one rank.
$ws.Range("A137").Value = "Birmania"
$ws.Range("B137").Value = 117
$ws.Range("C137").Value = 6
$ws.Range("D137").Value = 7
$ws.Range("E137").Value = 105
$ws.Range("H137").Value = 5

$ws.Range("A138").Value = "Trinidad yTobago"
$ws.Range("B138").Value = 114
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 21
$ws.Range("E138").Value = 85
$ws.Range("H138").Value = 8

# --- Guinea Ecuatorial (row 149) ---
$ws.Range("D149").Value = 7
$ws.Range("E149").Value = 72
